$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 107.928617
$ws.Range("H2").Value = 323.785851
$ws.Range("I2").Value = 0.2068777607879145
$ws.Range("J2").Value = 0.2068777607879145
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.701354
$ws.Range("N2").Value = 8.104061999999999
$ws.Range("O2").Value = 0.02221077311549548
$ws.Range("P2").Value = 0.02221077311549548
$ws.Range("Q2").Value = 291.5534012474179
$ws.Range("R2").Value = 2623.980611226762
$ws.Range("S2").Value = 0.004594915007502116
$ws.Range("T2").Value = 0.004594915007502116

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 107.928617
$ws.Range("H3").Value = 323.785851
$ws.Range("I3").Value = 0.2068777607879145
$ws.Range("J3").Value = 0.2068777607879145
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 88.14978533333333
$ws.Range("N3").Value = 264.449356
$ws.Range("O3").Value = 0.7247753838328104
$ws.Range("P3").Value = 0.7247753838328105
$ws.Range("Q3").Value = 9513.88441987355
$ws.Range("R3").Value = 85624.95977886194
$ws.Range("S3").Value = 0.149939908481533
$ws.Range("T3").Value = 0.1499399084815331

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 107.928617
$ws.Range("H4").Value = 323.785851
$ws.Range("I4").Value = 0.2068777607879145
$ws.Range("J4").Value = 0.2068777607879145
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.24063
$ws.Range("N4").Value = 0.72189
$ws.Range("O4").Value = 0.001978481285600361
$ws.Range("P4").Value = 0.001978481285600361
$ws.Range("Q4").Value = 25.97086310871
$ws.Range("R4").Value = 233.73776797839
$ws.Range("S4").Value = 0.0004093037781257971
$ws.Range("T4").Value = 0.0004093037781257971

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 107.928617
$ws.Range("H5").Value = 323.785851
$ws.Range("I5").Value = 0.2068777607879145
$ws.Range("J5").Value = 0.2068777607879145
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 30.53182233333333
$ws.Range("N5").Value = 91.595467
$ws.Range("O5").Value = 0.2510353617660938
$ws.Range("P5").Value = 0.2510353617660938
$ws.Range("Q5").Value = 3295.257358926379
$ws.Range("R5").Value = 29657.31623033742
$ws.Range("S5").Value = 0.05193363352075353
$ws.Range("T5").Value = 0.05193363352075353

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 143.300008
$ws.Range("H6").Value = 429.900024
$ws.Range("I6").Value = 0.2746777045788536
$ws.Range("J6").Value = 0.2746777045788536
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.701354
$ws.Range("N6").Value = 8.104061999999999
$ws.Range("O6").Value = 0.02221077311549548
$ws.Range("P6").Value = 0.02221077311549548
$ws.Range("Q6").Value = 387.104049810832
$ws.Range("R6").Value = 3483.936448297488
$ws.Range("S6").Value = 0.006100804176286011
$ws.Range("T6").Value = 0.006100804176286011

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 143.300008
$ws.Range("H7").Value = 429.900024
$ws.Range("I7").Value = 0.2746777045788536
$ws.Range("J7").Value = 0.2746777045788536
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 88.14978533333333
$ws.Range("N7").Value = 264.449356
$ws.Range("O7").Value = 0.7247753838328104
$ws.Range("P7").Value = 0.7247753838328105
$ws.Range("Q7").Value = 12631.86494346495
$ws.Range("R7").Value = 113686.7844911845
$ws.Range("S7").Value = 0.1990796387664539
$ws.Range("T7").Value = 0.1990796387664539

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 143.300008
$ws.Range("H8").Value = 429.900024
$ws.Range("I8").Value = 0.2746777045788536
$ws.Range("J8").Value = 0.2746777045788536
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.24063
$ws.Range("N8").Value = 0.72189
$ws.Range("O8").Value = 0.001978481285600361
$ws.Range("P8").Value = 0.001978481285600361
$ws.Range("Q8").Value = 34.48228092504001
$ws.Range("R8").Value = 310.34052832536
$ws.Range("S8").Value = 0.0005434446980809265
$ws.Range("T8").Value = 0.0005434446980809265

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 143.300008
$ws.Range("H9").Value = 429.900024
$ws.Range("I9").Value = 0.2746777045788536
$ws.Range("J9").Value = 0.2746777045788536
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 30.53182233333333
$ws.Range("N9").Value = 91.595467
$ws.Range("O9").Value = 0.2510353617660938
$ws.Range("P9").Value = 0.2510353617660938
$ws.Range("Q9").Value = 4375.210384621246
$ws.Range("R9").Value = 39376.89346159121
$ws.Range("S9").Value = 0.06895381693803275
$ws.Range("T9").Value = 0.06895381693803275

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 134.5459086666667
$ws.Range("H10").Value = 403.637726
$ws.Range("I10").Value = 0.2578978317505474
$ws.Range("J10").Value = 0.2578978317505473
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.701354
$ws.Range("N10").Value = 8.104061999999999
$ws.Range("O10").Value = 0.02221077311549548
$ws.Range("P10").Value = 0.02221077311549548
$ws.Range("Q10").Value = 363.4561285603347
$ws.Range("R10").Value = 3271.105157043012
$ws.Range("S10").Value = 0.005728110227989634
$ws.Range("T10").Value = 0.005728110227989633

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 134.5459086666667
$ws.Range("H11").Value = 403.637726
$ws.Range("I11").Value = 0.2578978317505474
$ws.Range("J11").Value = 0.2578978317505473
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 88.14978533333333
$ws.Range("N11").Value = 264.449356
$ws.Range("O11").Value = 0.7247753838328104
$ws.Range("P11").Value = 0.7247753838328105
$ws.Range("Q11").Value = 11860.19296644494
$ws.Range("R11").Value = 106741.7366980045
$ws.Range("S11").Value = 0.1869179999966525
$ws.Range("T11").Value = 0.1869179999966525

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 134.5459086666667
$ws.Range("H12").Value = 403.637726
$ws.Range("I12").Value = 0.2578978317505474
$ws.Range("J12").Value = 0.2578978317505473
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.24063
$ws.Range("N12").Value = 0.72189
$ws.Range("O12").Value = 0.001978481285600361
$ws.Range("P12").Value = 0.001978481285600361
$ws.Range("Q12").Value = 32.37578200246001
$ws.Range("R12").Value = 291.38203802214
$ws.Range("S12").Value = 0.0005102460337153686
$ws.Range("T12").Value = 0.0005102460337153685

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 134.5459086666667
$ws.Range("H13").Value = 403.637726
$ws.Range("I13").Value = 0.2578978317505474
$ws.Range("J13").Value = 0.2578978317505473
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 30.53182233333333
$ws.Range("N13").Value = 91.595467
$ws.Range("O13").Value = 0.2510353617660938
$ws.Range("P13").Value = 0.2510353617660938
$ws.Range("Q13").Value = 4107.931779087561
$ws.Range("R13").Value = 36971.38601178805
$ws.Range("S13").Value = 0.06474147549218985
$ws.Range("T13").Value = 0.06474147549218984

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 135.927831
$ws.Range("H14").Value = 407.783493
$ws.Range("I14").Value = 0.2605467028826847
$ws.Range("J14").Value = 0.2605467028826847
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 2.701354
$ws.Range("N14").Value = 8.104061999999999
$ws.Range("O14").Value = 0.02221077311549548
$ws.Range("P14").Value = 0.02221077311549548
$ws.Range("Q14").Value = 367.1891899831739
$ws.Range("R14").Value = 3304.702709848566
$ws.Range("S14").Value = 0.005786943703717721
$ws.Range("T14").Value = 0.005786943703717721

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 135.927831
$ws.Range("H15").Value = 407.783493
$ws.Range("I15").Value = 0.2605467028826847
$ws.Range("J15").Value = 0.2605467028826847
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 88.14978533333333
$ws.Range("N15").Value = 264.449356
$ws.Range("O15").Value = 0.7247753838328104
$ws.Range("P15").Value = 0.7247753838328105
$ws.Range("Q15").Value = 11982.00912347561
$ws.Range("R15").Value = 107838.0821112805
$ws.Range("S15").Value = 0.188837836588171
$ws.Range("T15").Value = 0.188837836588171

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 135.927831
$ws.Range("H16").Value = 407.783493
$ws.Range("I16").Value = 0.2605467028826847
$ws.Range("J16").Value = 0.2605467028826847
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.24063
$ws.Range("N16").Value = 0.72189
$ws.Range("O16").Value = 0.001978481285600361
$ws.Range("P16").Value = 0.001978481285600361
$ws.Range("Q16").Value = 32.70831397353
$ws.Range("R16").Value = 294.37482576177
$ws.Range("S16").Value = 0.0005154867756782693
$ws.Range("T16").Value = 0.0005154867756782693

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 135.927831
$ws.Range("H17").Value = 407.783493
$ws.Range("I17").Value = 0.2605467028826847
$ws.Range("J17").Value = 0.2605467028826847
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 30.53182233333333
$ws.Range("N17").Value = 91.595467
$ws.Range("O17").Value = 0.2510353617660938
$ws.Range("P17").Value = 0.2510353617660938
$ws.Range("Q17").Value = 4150.124386247359
$ws.Range("R17").Value = 37351.11947622623
$ws.Range("S17").Value = 0.06540643581511771
$ws.Range("T17").Value = 0.06540643581511771
